$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Wood"

# Update reference text strings (cells B26:B29) with the expanded references.
# Order matters here for shared-string table append order (matches upstream edit history).
$ws.Range("B26").Value = "[1] Eurostat (2011) Forestry in the EU and the world - a statistical portrait (http://refman.et-model.com/publications/1877)"
$ws.Range("B28").Value = "[3] OECD/IEA (2005) Energy statistics manual (http://refman.et-model.com/publications/1094)"
$ws.Range("B29").Value = "[4] Eurostat (2012) Roundwood, fuelwood and other basic products (for_basic) (http://appsso.eurostat.ec.europa.eu/nui/show.do?dataset=for_basic&lang=en)"
$ws.Range("B27").Value = "[2] IEA (2014) Unit Converter (http://www.iea.org/statistics/resources/unitconverter/)"

# Update the active cell selection to match the new state
[void]$ws.Range("B28").Select()
